$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell edits (order matters for shared-string pool ordering) ---

# Row 19 - Nombre lista cliente - new "Efigas" column entry
$ws.Range("E19").Value = "Efigas"

# Rows 8-13 - new plain values (no prior style) for Efigas convenio
$ws.Range("E8").Value = "30"
$ws.Range("E9").Value = "35"
$ws.Range("E10").Value = "44"
$ws.Range("E11").Value = "47"
$ws.Range("E12").Value = "54"
$ws.Range("E13").Value = "7709998002425"

# Rows 14, 16, 17 - brand new cells that need Text number format applied
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "192.141.2.1"

$ws.Range("E15").Value = "65"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "Efigas/recuados"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "administrador"

# Row 3 - nombre convenio
$ws.Range("E3").Value = "efigas"

# Rows reusing already-existing shared strings
$ws.Range("E5").Value = "4"
$ws.Range("E6").Value = "16"
$ws.Range("E7").Value = "21"
$ws.Range("E18").Value = "123"

# --- View/selection changes ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E13").Select()

$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
